$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-02-18 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-02-19 Wednesday", 2) | Out-Null
$d.Content.Find.Execute("59÷5=11, 4", $true, $false, $false, $false, $false, $true, 1, $false, "76÷3=25, 1", 2) | Out-Null
$d.Content.Find.Execute("88÷2=44, 0", $true, $false, $false, $false, $false, $true, 1, $false, "89÷5=17, 4", 2) | Out-Null
$d.Content.Find.Execute("52÷5=10, 2", $true, $false, $false, $false, $false, $true, 1, $false, "34÷2=17, 0", 2) | Out-Null
$d.Content.Find.Execute("46÷7=6, 4", $true, $false, $false, $false, $false, $true, 1, $false, "11÷7=1, 4", 2) | Out-Null
$d.Content.Find.Execute("29÷8=3, 5", $true, $false, $false, $false, $false, $true, 1, $false, "24÷6=4, 0", 2) | Out-Null
$d.Content.Find.Execute("44÷9=4, 8", $true, $false, $false, $false, $false, $true, 1, $false, "95÷5=19, 0", 2) | Out-Null
$d.Content.Find.Execute("99÷9=11, 0", $true, $false, $false, $false, $false, $true, 1, $false, "43÷5=8, 3", 2) | Out-Null
$d.Content.Find.Execute("56÷8=7, 0", $true, $false, $false, $false, $false, $true, 1, $false, "90÷8=11, 2", 2) | Out-Null
$d.Content.Find.Execute("53÷7=7, 4", $true, $false, $false, $false, $false, $true, 1, $false, "11÷2=5, 1", 2) | Out-Null
$d.Content.Find.Execute("19÷6=3, 1", $true, $false, $false, $false, $false, $true, 1, $false, "40÷5=8, 0", 2) | Out-Null
$d.Content.Find.Execute("96÷6=16, 0", $true, $false, $false, $false, $false, $true, 1, $false, "98÷7=14, 0", 2) | Out-Null
$d.Content.Find.Execute("55÷5=11, 0", $true, $false, $false, $false, $false, $true, 1, $false, "12÷8=1, 4", 2) | Out-Null
$d.Content.Find.Execute("82÷2=41, 0", $true, $false, $false, $false, $false, $true, 1, $false, "54÷8=6, 6", 2) | Out-Null
$d.Content.Find.Execute("71÷3=23, 2", $true, $false, $false, $false, $false, $true, 1, $false, "16÷9=1, 7", 2) | Out-Null
$d.Content.Find.Execute("13÷7=1, 6", $true, $false, $false, $false, $false, $true, 1, $false, "43÷4=10, 3", 2) | Out-Null
$d.Content.Find.Execute("62÷4=15, 2", $true, $false, $false, $false, $false, $true, 1, $false, "73÷8=9, 1", 2) | Out-Null
$d.Content.Find.Execute("36÷4=9, 0", $true, $false, $false, $false, $false, $true, 1, $false, "19÷7=2, 5", 2) | Out-Null
$d.Content.Find.Execute("28÷3=9, 1", $true, $false, $false, $false, $false, $true, 1, $false, "43÷7=6, 1", 2) | Out-Null
$d.Content.Find.Execute("96÷5=19, 1", $true, $false, $false, $false, $false, $true, 1, $false, "16÷4=4, 0", 2) | Out-Null
$d.Content.Find.Execute("58÷4=14, 2", $true, $false, $false, $false, $false, $true, 1, $false, "90÷8=11, 2", 2) | Out-Null
$d.Content.Find.Execute("48÷4=12, 0", $true, $false, $false, $false, $false, $true, 1, $false, "53÷7=7, 4", 2) | Out-Null
$d.Content.Find.Execute("89÷7=12, 5", $true, $false, $false, $false, $false, $true, 1, $false, "29÷9=3, 2", 2) | Out-Null
$d.Content.Find.Execute("54÷6=9, 0", $true, $false, $false, $false, $false, $true, 1, $false, "55÷6=9, 1", 2) | Out-Null
$d.Content.Find.Execute("29÷2=14, 1", $true, $false, $false, $false, $false, $true, 1, $false, "29÷9=3, 2", 2) | Out-Null
$d.Content.Find.Execute("33÷7=4, 5", $true, $false, $false, $false, $false, $true, 1, $false, "52÷5=10, 2", 2) | Out-Null
